$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3 describe the same sighting but the "Id" (A), "Antal"/"Enhet"
# (I/J) and coordinate/accuracy (Q/R/S) values were swapped between the two
# rows. Swap them back using a scratch cell far outside the used range so
# cell types (number vs. text) are preserved exactly via Range.Copy instead
# of re-typing raw values (which would coerce "10" into the number 10).
$cols = @("A", "I", "J", "Q", "R", "S")
$scratchRow = 1000

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $scratch = $ws.Range($col + $scratchRow)

    $scratch.ClearContents()
    $cell2.Copy($scratch)    # scratch now holds the old row-2 value
    $cell2.ClearContents()
    $cell3.Copy($cell2)      # row 2 now holds the old row-3 value
    $cell3.ClearContents()
    $scratch.Copy($cell3)    # row 3 now holds the old row-2 value
    $scratch.ClearContents()
}
